$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.424537301063538
$ws.Range("B1").Value = 3.460838317871094
$ws.Range("C1").Value = 5.436500072479248
$ws.Range("D1").Value = 1.742715001106262
$ws.Range("E1").Value = 0.9777151942253113
